$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B4 sales value
$ws.Range("B4").Value = 3750000

# Convert row 2 totals into a shared formula across B2:M2
$ws.Range("B2:M2").Formula = "=SUM(B3:B7)"

# Restore the active selection to B4 as left by the author
$ws.Range("B4").Select() | Out-Null
